# Apply ASR Results data update (M09 Data Aug Froze D3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-18, columns A (word), B (word), C (count)
$data = @(
    @("<made>",    "<may>",    8),
    @("<about>",   "<about>",  11),
    @("<have>",    "<have>",   8),
    @("<an>",      "<an>",     14),
    @("<nine>",    "<nine>",   16),
    @("<left>",    "<left>",   6),
    @("<oscar>",   "<oscar>",  11),
    @("<line>",    "<line>",   13),
    @("<alt>",     "<up>",     12),
    @("<see>",     "<see>",    9),
    @("<seven>",   "<seven>",  10),
    @("<command>", "<command>",13),
    @("<day>",     "<day>",    7),
    @("<victor>",  "<echo>",   9),
    @("<we>",      "<we>",     11),
    @("<for>",     "<four>",   11),
    @("<out>",     "<out>",    10)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
